$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (the only semantic change per the commit message).
$ws.Name = "variable_mapping"

# The target column widths (in OOXML "character" width units) from the diff.
# The COM ColumnWidth setter adds Excel's standard 5px (=5/6 character-unit)
# padding on save, so we back that padding out before assigning so the
# round-tripped stored width lands as close as possible to the target.
$targetWidths = @(14.48, 11.25, 10.21, 14.48, 13.33, 19.79, 16.56, 18.75, 13.33, 14.48, 20.83, 15.52, 27.19, 15.52, 10.21, 12.29)
$padding = 5.0 / 6.0

for ($i = 0; $i -lt $targetWidths.Length; $i++) {
    $col = $i + 1
    $ws.Columns.Item($col).ColumnWidth = $targetWidths[$i] - $padding
}
